# Generate Report for Handback
# Updates the localization-status workbook to reflect that the
# "20bcb66c-dd1a-43b3-a248-76b294441b45" file has now been handed back
# (in sync with en-US), adding the Latest Target File / Latest Handback
# File references and the Latest Handback DateTime for both the zh-cn
# and de-de language sheets, and refreshing the status shown on the
# Overview sheet.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------
# Overview sheet: the 20bcb66c row (row 2) status changes for both
# the zh-cn and de-de status columns.
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus

# ---------------------------------------------------------------
# Helper data describing the per-language hyperlink targets so the
# zh-cn and de-de sheets can be updated with the same logic.
# ---------------------------------------------------------------
$langs = @(
    @{
        SheetName = "zh-cn"
        MdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/6e684c28b587a70f0cfd5f5c6ff1ff3adb27a702/e2e/20bcb66c-dd1a-43b3-a248-76b294441b45.md"
        MdDisplay = "20bcb66c-dd1a-43b3-a248-76b294441b45.md"
        XlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/579fedc10ff59a7149311f4f12d601c1d77743ec/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/20bcb66c-dd1a-43b3-a248-76b294441b45.633f0c36d9a87b0eefdf682d49f1e52b4ef346d6.zh-cn.xlf"
        XlfDisplay = "20bcb66c-dd1a-43b3-a248-76b294441b45.633f0c36d9a87b0eefdf682d49f1e52b4ef346d6.zh-cn.xlf"
        Md2Target = "https://github.com/OpenLocalizationTest/oltest/blob/4a68d5de10ad6c3300541df7be93d85829774398/e2e/b68a2137-e2a6-439b-88bb-be13c4dbe045.md"
        Md2Display = "b68a2137-e2a6-439b-88bb-be13c4dbe045.md"
        Xlf2Target = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7cf28389d94558a735df52069639c0c7a33fa611/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b68a2137-e2a6-439b-88bb-be13c4dbe045.cafe1cadf14af19aedb092ec9298155126eb95ed.zh-cn.xlf"
        Xlf2Display = "b68a2137-e2a6-439b-88bb-be13c4dbe045.cafe1cadf14af19aedb092ec9298155126eb95ed.zh-cn.xlf"
        HandbackDateTime = "2016-03-24 04:38:57"
    },
    @{
        SheetName = "de-de"
        MdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/6e684c28b587a70f0cfd5f5c6ff1ff3adb27a702/e2e/20bcb66c-dd1a-43b3-a248-76b294441b45.md"
        MdDisplay = "20bcb66c-dd1a-43b3-a248-76b294441b45.md"
        XlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/32900312ec1e6af4822ada052026ac7daaba561d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/20bcb66c-dd1a-43b3-a248-76b294441b45.633f0c36d9a87b0eefdf682d49f1e52b4ef346d6.de-de.xlf"
        XlfDisplay = "20bcb66c-dd1a-43b3-a248-76b294441b45.633f0c36d9a87b0eefdf682d49f1e52b4ef346d6.de-de.xlf"
        Md2Target = "https://github.com/OpenLocalizationTest/oltest/blob/4a68d5de10ad6c3300541df7be93d85829774398/e2e/b68a2137-e2a6-439b-88bb-be13c4dbe045.md"
        Md2Display = "b68a2137-e2a6-439b-88bb-be13c4dbe045.md"
        Xlf2Target = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/52aaa75309f84d13edc3e411a90d4758a4ff3139/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b68a2137-e2a6-439b-88bb-be13c4dbe045.cafe1cadf14af19aedb092ec9298155126eb95ed.de-de.xlf"
        Xlf2Display = "b68a2137-e2a6-439b-88bb-be13c4dbe045.cafe1cadf14af19aedb092ec9298155126eb95ed.de-de.xlf"
        HandbackDateTime = "2016-03-24 04:39:04"
    }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.SheetName)

    # Status column (C) for the 20bcb66c row is now "handed back".
    $ws.Range("C2").Value = $newStatus

    # Latest Handback DateTime (H) for the 20bcb66c row.
    $ws.Range("H2").Value = $lang.HandbackDateTime

    # Rebuild the hyperlinks collection so that the two new links
    # (Latest Target File / Latest Handback File) land between the
    # existing row 2 and row 3 links, matching how the relationships
    # are ordered/renumbered.
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $lang.MdTarget, "", "", $lang.MdDisplay)
    $ws.Hyperlinks.Add($ws.Range("D2"), $lang.XlfTarget, "", "", $lang.XlfDisplay)

    $ws.Hyperlinks.Add($ws.Range("F2"), $lang.MdTarget, "", "", $lang.MdDisplay)
    $ws.Range("F2").Font.Underline = $true
    $ws.Range("F2").Font.Color = 15570276

    $ws.Hyperlinks.Add($ws.Range("G2"), $lang.XlfTarget, "", "", $lang.XlfDisplay)
    $ws.Range("G2").Font.Underline = $true
    $ws.Range("G2").Font.Color = 15570276

    $ws.Hyperlinks.Add($ws.Range("A3"), $lang.Md2Target, "", "", $lang.Md2Display)
    $ws.Hyperlinks.Add($ws.Range("D3"), $lang.Xlf2Target, "", "", $lang.Xlf2Display)
}
